$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")

# Add the new "Incluido" column header, matching the style of A1
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B1").Value = "Incluido"

# Mark rows 2-14 (the included users) with an "X" in column B
$ws.Range("B2:B14").Value = "X"

# Update the selection to match the saved state
$excel.Goto($ws.Range("A15:A18"))
